# Updated cryptos list values to match latest scrape (price + 1h volume % changes,
# plus the row-shift in the coin ranking list for rows 30-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking "Price" strings to stay text (matches the source data,
# which stores prices like "1.00" / "0.0000178" as literal text, not numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.198.88'
$ws.Range('E2').Value = '  +3.01%  '
$ws.Range('D3').Value = '2.408.74'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '571.87'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '144.59'
$ws.Range('E6').Value = '  +4.63%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '2.430.22'
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('E10').Value = '  +4.59%  '
$ws.Range('E11').Value = '  +0.89%  '
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('E13').Value = '  +3.47%  '
$ws.Range('D14').Value = '26.76'
$ws.Range('E14').Value = '  +4.34%  '
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +7.70%  '
$ws.Range('D16').Value = '2.850.21'
$ws.Range('E16').Value = '  +0.93%  '
$ws.Range('D17').Value = '62.088.84'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '2.433.30'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('E19').Value = '  -4.93%  '
$ws.Range('E20').Value = '  +2.53%  '
$ws.Range('D21').Value = '325.32'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('E23').Value = '  +14.02%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '65.32'
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').Value = '611.88'
$ws.Range('E26').Value = '  +10.56%  '
$ws.Range('D27').Value = '8.43'
$ws.Range('E27').Value = '  +5.35%  '
$ws.Range('D28').Value = '0.0₃0980'
$ws.Range('E28').Value = '  +8.07%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '8.07'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.40'
$ws.Range('E31').Value = '  +8.46%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.138'
$ws.Range('E32').Value = '  +7.08%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.84'
$ws.Range('E33').Value = '  +2.22%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  +5.11%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.996'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '4.75'
$ws.Range('E36').Value = '  +5.44%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '152.79'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '0.372'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '5.37'
$ws.Range('E39').Value = '  +6.34%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.51'
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.73'
$ws.Range('E41').Value = '  +19.23%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.75'
$ws.Range('E42').Value = '  +7.03%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '42.29'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0280'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '143.66'
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '3.57'
$ws.Range('E47').Value = '  +2.70%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '20.19'
$ws.Range('E48').Value = '  +7.09%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.598'
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D50').Value = '0.0513'
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.0913'
$ws.Range('E51').Value = '  +1.76%  '
